{"js": "// Update master to output generated at 4250d90\n// Updates the date heading and the 25 \"NN\u00d7NN=\" practice-table equations.\n//\n// Both the heading paragraph and every table cell are targeted by their\n// fixed position (first paragraph; table 1, row/col) rather than by a\n// document-wide text search. Several of the old/new equation strings are\n// not unique across the document (e.g. \"71\u00d787=\" is the *old* text of the\n// cell at row 0/col 3 but also the *new* text written into the cell at\n// row 9/col 0, and vice versa for \"11\u00d792=\"), so replacing by searching\n// for matching text globally would risk rewriting the wrong cell, or\n// re-matching a cell that was already updated earlier in the same run.\n// Addressing cells by their fixed position avoids all of that ambiguity.\n\n// --- Date heading -------------------------------------------------------\nconst firstPara = context.document.body.paragraphs.getFirst();\nfirstPara.load(\"text\");\n\n// --- Equation table -------------------------------------------------------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\n\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst dateOld = \"2024-06-19 Wednesday\";\nconst dateNew = \"2024-06-20 Thursday\";\n\nconst cellEdits = [\n  { row: 0, col: 0, oldValue: \"87\u00d765=\", newValue: \"69\u00d762=\" },\n  { row: 0, col: 1, oldValue: \"86\u00d784=\", newValue: \"85\u00d716=\" },\n  { row: 0, col: 2, oldValue: \"37\u00d740=\", newValue: \"34\u00d796=\" },\n  { row: 0, col: 3, oldValue: \"92\u00d743=\", newValue: \"71\u00d787=\" },\n  { row: 0, col: 4, oldValue: \"84\u00d778=\", newValue: \"72\u00d779=\" },\n  { row: 4, col: 0, oldValue: \"50\u00d739=\", newValue: \"86\u00d756=\" },\n  { row: 4, col: 1, oldValue: \"68\u00d789=\", newValue: \"64\u00d716=\" },\n  { row: 4, col: 2, oldValue: \"67\u00d767=\", newValue: \"48\u00d719=\" },\n  { row: 4, col: 3, oldValue: \"13\u00d789=\", newValue: \"52\u00d773=\" },\n  { row: 4, col: 4, oldValue: \"96\u00d799=\", newValue: \"70\u00d745=\" },\n  { row: 9, col: 0, oldValue: \"71\u00d787=\", newValue: \"11\u00d792=\" },\n  { row: 9, col: 1, oldValue: \"51\u00d748=\", newValue: \"23\u00d793=\" },\n  { row: 9, col: 2, oldValue: \"73\u00d764=\", newValue: \"93\u00d727=\" },\n  { row: 9, col: 3, oldValue: \"89\u00d798=\", newValue: \"70\u00d753=\" },\n  { row: 9, col: 4, oldValue: \"53\u00d797=\", newValue: \"97\u00d795=\" },\n  { row: 14, col: 0, oldValue: \"87\u00d737=\", newValue: \"97\u00d789=\" },\n  { row: 14, col: 1, oldValue: \"29\u00d718=\", newValue: \"37\u00d795=\" },\n  { row: 14, col: 2, oldValue: \"82\u00d769=\", newValue: \"97\u00d773=\" },\n  { row: 14, col: 3, oldValue: \"81\u00d751=\", newValue: \"16\u00d731=\" },\n  { row: 14, col: 4, oldValue: \"47\u00d775=\", newValue: \"82\u00d777=\" },\n  { row: 19, col: 0, oldValue: \"38\u00d719=\", newValue: \"33\u00d741=\" },\n  { row: 19, col: 1, oldValue: \"64\u00d721=\", newValue: \"29\u00d715=\" },\n  { row: 19, col: 2, oldValue: \"84\u00d779=\", newValue: \"45\u00d715=\" },\n  { row: 19, col: 3, oldValue: \"76\u00d735=\", newValue: \"62\u00d747=\" },\n  { row: 19, col: 4, oldValue: \"26\u00d783=\", newValue: \"50\u00d738=\" }\n];\n\nconst cells = cellEdits.map((edit) => table.getCell(edit.row, edit.col));\ncells.forEach((cell) => cell.load(\"value\"));\n\nawait context.sync();\n\nif (firstPara.text.trim() === dateOld) {\n  firstPara.getRange().insertText(dateNew, Word.InsertLocation.replace);\n}\n\ncellEdits.forEach((edit, i) => {\n  const cell = cells[i];\n  if (cell.value === edit.oldValue) {\n    cell.value = edit.newValue;\n  }\n});\n\nawait context.sync();\n", "ps1": "# Update master to output generated at 4250d90\n# Updates the date heading and the 25 \"NN\u00d7NN=\" practice-table equations.\n#\n# Both the heading paragraph and every table cell are targeted by their\n# fixed position (paragraph 1; table 1, row/col) rather than by searching\n# for the old text globally. Several of the old/new equation strings are\n# not unique within the document (e.g. \"71\u00d787=\" is the *old* text of one\n# cell and also the *new* text written into a different cell), so a\n# document-wide Find/Replace-All would corrupt cells that were already\n# updated (or not-yet-updated). Positional addressing avoids that.\n\n$d = $word.ActiveDocument\n\n# --- Date heading -----------------------------------------------------\n$dateOld = '2024-06-19 Wednesday'\n$dateNew = '2024-06-20 Thursday'\n\n$headingPara = $d.Paragraphs.Item(1)\nif ($headingPara.Range.Text.TrimEnd([char]13, [char]7) -eq $dateOld) {\n    $headingPara.Range.Text = $dateNew\n}\n\n# --- Equation table -----------------------------------------------------\n$table = $d.Tables.Item(1)\n\n$cellEdits = @(\n    @{ Row = 1; Col = 1; Old = '87\u00d765='; New = '69\u00d762=' },\n    @{ Row = 1; Col = 2; Old = '86\u00d784='; New = '85\u00d716=' },\n    @{ Row = 1; Col = 3; Old = '37\u00d740='; New = '34\u00d796=' },\n    @{ Row = 1; Col = 4; Old = '92\u00d743='; New = '71\u00d787=' },\n    @{ Row = 1; Col = 5; Old = '84\u00d778='; New = '72\u00d779=' },\n    @{ Row = 5; Col = 1; Old = '50\u00d739='; New = '86\u00d756=' },\n    @{ Row = 5; Col = 2; Old = '68\u00d789='; New = '64\u00d716=' },\n    @{ Row = 5; Col = 3; Old = '67\u00d767='; New = '48\u00d719=' },\n    @{ Row = 5; Col = 4; Old = '13\u00d789='; New = '52\u00d773=' },\n    @{ Row = 5; Col = 5; Old = '96\u00d799='; New = '70\u00d745=' },\n    @{ Row = 10; Col = 1; Old = '71\u00d787='; New = '11\u00d792=' },\n    @{ Row = 10; Col = 2; Old = '51\u00d748='; New = '23\u00d793=' },\n    @{ Row = 10; Col = 3; Old = '73\u00d764='; New = '93\u00d727=' },\n    @{ Row = 10; Col = 4; Old = '89\u00d798='; New = '70\u00d753=' },\n    @{ Row = 10; Col = 5; Old = '53\u00d797='; New = '97\u00d795=' },\n    @{ Row = 15; Col = 1; Old = '87\u00d737='; New = '97\u00d789=' },\n    @{ Row = 15; Col = 2; Old = '29\u00d718='; New = '37\u00d795=' },\n    @{ Row = 15; Col = 3; Old = '82\u00d769='; New = '97\u00d773=' },\n    @{ Row = 15; Col = 4; Old = '81\u00d751='; New = '16\u00d731=' },\n    @{ Row = 15; Col = 5; Old = '47\u00d775='; New = '82\u00d777=' },\n    @{ Row = 20; Col = 1; Old = '38\u00d719='; New = '33\u00d741=' },\n    @{ Row = 20; Col = 2; Old = '64\u00d721='; New = '29\u00d715=' },\n    @{ Row = 20; Col = 3; Old = '84\u00d779='; New = '45\u00d715=' },\n    @{ Row = 20; Col = 4; Old = '76\u00d735='; New = '62\u00d747=' },\n    @{ Row = 20; Col = 5; Old = '26\u00d783='; New = '50\u00d738=' }\n)\n\nforeach ($edit in $cellEdits) {\n    $cell = $table.Cell($edit.Row, $edit.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -eq $edit.Old) {\n        $cell.Range.Text = $edit.New\n    }\n}\n"}
